$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the "Brief:" title
# paragraph. It needs to move to the end of the new MVP bullet item, so
# remove it from its current spot first.
$d.Bookmarks.Item("_GoBack").Delete()

# --- Add the "MVP:" heading paragraph at the end of the document ---
$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc - 1, $endOfDoc - 1)
$insertionPoint.InsertParagraphAfter()

$mvpRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$mvpRange.Text = "MVP:"
$mvpPara = $d.Paragraphs.Last
$mvpPara.Style = "Title"

# --- Add the bulleted "Need Start screen, game loop," paragraph ---
$endOfDoc2 = $d.Content.End
$insertionPoint2 = $d.Range($endOfDoc2 - 1, $endOfDoc2 - 1)
$insertionPoint2.InsertParagraphAfter()

# Append a trailing sentinel character so the range's end does not land
# exactly on the paragraph mark while we bookmark it (placing a bookmark
# collapsed right at a paragraph-mark offset gets mishandled); the
# sentinel is stripped again immediately after the bookmark is placed.
$listRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$listRange.Text = "Need Start screen, game loop,X"
$listPara = $d.Paragraphs.Last
$listPara.Style = "List Paragraph"
$listRange.ListFormat.ApplyBulletDefault()

$bookmarkRange = $d.Range($listRange.End - 1, $listRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$sentinel = $d.Range($listRange.End - 1, $listRange.End)
$sentinel.Delete()

# Line up the auto-generated "List Paragraph" style's priority with the
# value Word normally assigns to it.
$listParagraphStyle = $d.Styles.Item("List Paragraph")
$listParagraphStyle.Priority = 34
